$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 31,12
$data[0,0] = 2
$data[0,1] = 372
$data[0,2] = 16
$data[0,3] = 104
$data[0,4] = 369
$data[0,5] = 2325
$data[0,6] = 84
$data[0,7] = 747
$data[0,8] = 1120
$data[0,9] = 438
$data[0,10] = 125
$data[0,11] = 129
$data[1,0] = 2
$data[1,1] = 355
$data[1,2] = 17
$data[1,3] = 95
$data[1,4] = 375
$data[1,5] = 2536
$data[1,6] = 82
$data[1,7] = 585
$data[1,8] = 1113
$data[1,9] = 151
$data[1,10] = 125
$data[1,11] = 123
$data[2,0] = 2
$data[2,1] = 350
$data[2,2] = 15
$data[2,3] = 95
$data[2,4] = 367
$data[2,5] = 2316
$data[2,6] = 77
$data[2,7] = 543
$data[2,8] = 1096
$data[2,9] = 142
$data[2,10] = 126
$data[2,11] = 125
$data[3,0] = 2
$data[3,1] = 344
$data[3,2] = 17
$data[3,3] = 99
$data[3,4] = 381
$data[3,5] = 2539
$data[3,6] = 80
$data[3,7] = 565
$data[3,8] = 1099
$data[3,9] = 125
$data[3,10] = 133
$data[3,11] = 121
$data[4,0] = 2
$data[4,1] = 347
$data[4,2] = 15
$data[4,3] = 93
$data[4,4] = 374
$data[4,5] = 2301
$data[4,6] = 80
$data[4,7] = 546
$data[4,8] = 1118
$data[4,9] = 126
$data[4,10] = 125
$data[4,11] = 124
$data[5,0] = 2
$data[5,1] = 349
$data[5,2] = 17
$data[5,3] = 101
$data[5,4] = 378
$data[5,5] = 2458
$data[5,6] = 76
$data[5,7] = 567
$data[5,8] = 1095
$data[5,9] = 126
$data[5,10] = 137
$data[5,11] = 130
$data[6,0] = 2
$data[6,1] = 342
$data[6,2] = 16
$data[6,3] = 92
$data[6,4] = 368
$data[6,5] = 2471
$data[6,6] = 79
$data[6,7] = 588
$data[6,8] = 1085
$data[6,9] = 127
$data[6,10] = 130
$data[6,11] = 126
$data[7,0] = 2
$data[7,1] = 343
$data[7,2] = 16
$data[7,3] = 89
$data[7,4] = 376
$data[7,5] = 2566
$data[7,6] = 77
$data[7,7] = 696
$data[7,8] = 1092
$data[7,9] = 131
$data[7,10] = 135
$data[7,11] = 129
$data[8,0] = 2
$data[8,1] = 362
$data[8,2] = 18
$data[8,3] = 81
$data[8,4] = 390
$data[8,5] = 2303
$data[8,6] = 80
$data[8,7] = 668
$data[8,8] = 1056
$data[8,9] = 131
$data[8,10] = 131
$data[8,11] = 128
$data[9,0] = 2
$data[9,1] = 361
$data[9,2] = 21
$data[9,3] = 83
$data[9,4] = 396
$data[9,5] = 2232
$data[9,6] = 80
$data[9,7] = 604
$data[9,8] = 1071
$data[9,9] = 131
$data[9,10] = 132
$data[9,11] = 138
$data[10,0] = 2
$data[10,1] = 381
$data[10,2] = 18
$data[10,3] = 85
$data[10,4] = 402
$data[10,5] = 2331
$data[10,6] = 83
$data[10,7] = 573
$data[10,8] = 1053
$data[10,9] = 138
$data[10,10] = 130
$data[10,11] = 133
$data[11,0] = 2
$data[11,1] = 373
$data[11,2] = 21
$data[11,3] = 80
$data[11,4] = 408
$data[11,5] = 2354
$data[11,6] = 85
$data[11,7] = 562
$data[11,8] = 1062
$data[11,9] = 128
$data[11,10] = 140
$data[11,11] = 137
$data[12,0] = 2
$data[12,1] = 392
$data[12,2] = 17
$data[12,3] = 84
$data[12,4] = 417
$data[12,5] = 2309
$data[12,6] = 80
$data[12,7] = 542
$data[12,8] = 1056
$data[12,9] = 135
$data[12,10] = 137
$data[12,11] = 131
$data[13,0] = 2
$data[13,1] = 388
$data[13,2] = 20
$data[13,3] = 83
$data[13,4] = 418
$data[13,5] = 2312
$data[13,6] = 80
$data[13,7] = 554
$data[13,8] = 1069
$data[13,9] = 140
$data[13,10] = 135
$data[13,11] = 132
$data[14,0] = 2
$data[14,1] = 373
$data[14,2] = 17
$data[14,3] = 84
$data[14,4] = 393
$data[14,5] = 2333
$data[14,6] = 85
$data[14,7] = 521
$data[14,8] = 1059
$data[14,9] = 137
$data[14,10] = 132
$data[14,11] = 137
$data[15,0] = 2
$data[15,1] = 371
$data[15,2] = 19
$data[15,3] = 89
$data[15,4] = 408
$data[15,5] = 2238
$data[15,6] = 82
$data[15,7] = 547
$data[15,8] = 1059
$data[15,9] = 136
$data[15,10] = 135
$data[15,11] = 135
$data[16,0] = 2
$data[16,1] = 360
$data[16,2] = 17
$data[16,3] = 82
$data[16,4] = 374
$data[16,5] = 2285
$data[16,6] = 80
$data[16,7] = 537
$data[16,8] = 1056
$data[16,9] = 135
$data[16,10] = 136
$data[16,11] = 134
$data[17,0] = 2
$data[17,1] = 355
$data[17,2] = 19
$data[17,3] = 85
$data[17,4] = 389
$data[17,5] = 2323
$data[17,6] = 78
$data[17,7] = 534
$data[17,8] = 1054
$data[17,9] = 137
$data[17,10] = 139
$data[17,11] = 137
$data[18,0] = 2
$data[18,1] = 366
$data[18,2] = 17
$data[18,3] = 85
$data[18,4] = 397
$data[18,5] = 2278
$data[18,6] = 83
$data[18,7] = 517
$data[18,8] = 1069
$data[18,9] = 137
$data[18,10] = 135
$data[18,11] = 134
$data[19,0] = 2
$data[19,1] = 368
$data[19,2] = 19
$data[19,3] = 86
$data[19,4] = 409
$data[19,5] = 2311
$data[19,6] = 82
$data[19,7] = 533
$data[19,8] = 1063
$data[19,9] = 140
$data[19,10] = 140
$data[19,11] = 136
$data[20,0] = 2
$data[20,1] = 367
$data[20,2] = 18
$data[20,3] = 84
$data[20,4] = 388
$data[20,5] = 2344
$data[20,6] = 79
$data[20,7] = 539
$data[20,8] = 1079
$data[20,9] = 130
$data[20,10] = 131
$data[20,11] = 125
$data[21,0] = 1
$data[21,1] = 403
$data[21,2] = 58
$data[21,3] = 135
$data[21,4] = 452
$data[21,5] = 2375
$data[21,6] = 131
$data[21,7] = 540
$data[21,8] = 1083
$data[21,9] = 126
$data[21,10] = 131
$data[21,11] = 126
$data[22,0] = 3
$data[22,1] = 313
$data[22,2] = 19
$data[22,3] = 80
$data[22,4] = 314
$data[22,5] = 3993
$data[22,6] = 85
$data[22,7] = 359
$data[22,8] = 1887
$data[22,9] = 49
$data[22,10] = 43
$data[22,11] = 43
$data[23,0] = 13
$data[23,1] = 339
$data[23,2] = 36
$data[23,3] = 118
$data[23,4] = 342
$data[23,5] = 2537
$data[23,6] = 116
$data[23,7] = 390
$data[23,8] = 1142
$data[23,9] = 53
$data[23,10] = 46
$data[23,11] = 42
$data[24,0] = 15
$data[24,1] = 2119
$data[24,2] = 132
$data[24,3] = 5274
$data[24,4] = 2088
$data[24,5] = 2646
$data[24,6] = 5205
$data[24,7] = 393
$data[24,8] = 1153
$data[24,9] = 44
$data[24,10] = 46
$data[24,11] = 45
$data[25,0] = 64
$data[25,1] = 817
$data[25,2] = 108
$data[25,3] = 200
$data[25,4] = 985
$data[25,5] = 3396
$data[25,6] = 196
$data[25,7] = 752
$data[25,8] = 2174
$data[25,9] = 99
$data[25,10] = 98
$data[25,11] = 97
$data[26,0] = 106
$data[26,1] = 8005
$data[26,2] = 466
$data[26,3] = 9255
$data[26,4] = 8179
$data[26,5] = 3615
$data[26,6] = 9212
$data[26,7] = 798
$data[26,8] = 2211
$data[26,9] = 132
$data[26,10] = 123
$data[26,11] = 122
$data[27,0] = 145
$data[27,1] = 1066
$data[27,2] = 368
$data[27,3] = 416
$data[27,4] = 894
$data[27,5] = 2573
$data[27,6] = 244
$data[27,7] = 859
$data[27,8] = 1372
$data[27,9] = 463
$data[27,10] = 371
$data[27,11] = 415
$data[28,0] = 24
$data[28,1] = 449
$data[28,2] = 26
$data[28,3] = 181
$data[28,4] = 439
$data[28,5] = 2420
$data[28,6] = 211
$data[28,7] = 393
$data[28,8] = 1024
$data[28,9] = 152
$data[28,10] = 149
$data[28,11] = 170
$data[29,0] = 151
$data[29,1] = 2532416
$data[29,2] = 429
$data[29,3] = 354
$data[29,4] = 2644224
$data[29,5] = 2288
$data[29,6] = 246
$data[29,7] = 560
$data[29,8] = 1220
$data[29,9] = 356
$data[29,10] = 327
$data[29,11] = 350
$data[30,0] = 138
$data[30,1] = 10435
$data[30,2] = 939
$data[30,3] = 2414
$data[30,4] = 10102
$data[30,5] = 934
$data[30,6] = 2225
$data[30,7] = 111
$data[30,8] = 109
$data[30,9] = 112
$data[30,10] = 111
$data[30,11] = 109

$ws.Range("B2:M32").Value = $data

$ws.Range("B2:M32").Select()
